$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1099.6364
$ws.Range("I28").Value = 372.57144
$ws.Range("J28").Value = 2372
$ws.Range("K28").Value = 372.57144
$ws.Range("L28").Value = 2372
$ws.Range("M28").Value = 112.42856
$ws.Range("N28").Value = -3342

$ws.Range("H88").Value = 14290643
$ws.Range("I88").Value = 100000000
$ws.Range("K88").Value = 100000000
$ws.Range("M88").Value = -99999594

$ws.Range("H91").Value = 14290643
$ws.Range("I91").Value = 100000000
$ws.Range("K91").Value = 100000000
$ws.Range("M91").Value = -99998596

$ws.Range("H111").Value = 5064
$ws.Range("J111").Value = 2339.4
$ws.Range("L111").Value = 7018.200000000001
$ws.Range("N111").Value = -13152.2

$ws.Range("H113").Value = 3988.75
$ws.Range("I113").Value = 3079.6
$ws.Range("K113").Value = 3079.6
$ws.Range("M113").Value = 174.4000000000001

$ws.Range("H131").Value = 1002964.5
$ws.Range("I131").Value = 1430892.4
$ws.Range("K131").Value = 4292677.199999999
$ws.Range("M131").Value = -4287637.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1223.0769
$ws.Range("I45").Value = 1117.6818
$ws.Range("K45").Value = 1117.6818
$ws.Range("M45").Value = -740.6818000000001

$ws.Range("H63").Value = 4165.6665
$ws.Range("I63").Value = 1498.5
$ws.Range("J63").Value = 9500
$ws.Range("K63").Value = 1498.5
$ws.Range("L63").Value = 9500
$ws.Range("M63").Value = -812.5
$ws.Range("N63").Value = -10872

$ws.Range("H66").Value = 4165.6665
$ws.Range("I66").Value = 1498.5
$ws.Range("J66").Value = 9500
$ws.Range("K66").Value = 7492.5
$ws.Range("L66").Value = 47500
$ws.Range("M66").Value = -4060.5
$ws.Range("N66").Value = -54364

$ws.Range("H74").Value = 148121.7
$ws.Range("I74").Value = 215197.23
$ws.Range("K74").Value = 215197.23
$ws.Range("M74").Value = -214323.23

$ws.Range("H77").Value = 148121.7
$ws.Range("I77").Value = 215197.23
$ws.Range("K77").Value = 1075986.15
$ws.Range("M77").Value = -1071618.15

$ws.Range("H102").Value = 2345.3044
$ws.Range("I102").Value = 1774.7222
$ws.Range("K102").Value = 1774.7222
$ws.Range("M102").Value = -152.7221999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2748650.5
$ws.Range("I107").Value = 3206485.8
$ws.Range("J107").Value = 1639.25
$ws.Range("K107").Value = 3206485.8
$ws.Range("L107").Value = 1639.25
$ws.Range("M107").Value = -3204565.8
$ws.Range("N107").Value = -5479.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3477517.5
$ws.Range("I31").Value = 3729.6428
$ws.Range("K31").Value = 3729.6428
$ws.Range("M31").Value = -3434.6428

$ws.Range("H34").Value = 3477517.5
$ws.Range("I34").Value = 3729.6428
$ws.Range("K34").Value = 3729.6428
$ws.Range("M34").Value = -3527.6428

$ws.Range("H132").Value = 3579.9556
$ws.Range("I132").Value = 3495.1853
$ws.Range("J132").Value = 3707.111
$ws.Range("K132").Value = 10485.5559
$ws.Range("L132").Value = 11121.333
$ws.Range("M132").Value = -7955.555899999999
$ws.Range("N132").Value = -16181.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H107").Value = 406.3846
$ws.Range("I107").Value = 34
$ws.Range("J107").Value = 437.41666
$ws.Range("K107").Value = 102
$ws.Range("L107").Value = 1312.24998
$ws.Range("M107").Value = 1818
$ws.Range("N107").Value = -5152.249980000001

$ws.Range("H122").Value = 1011.75
$ws.Range("J122").Value = 2473
$ws.Range("L122").Value = 22257
$ws.Range("N122").Value = -27157

$ws.Range("H129").Value = 1607.5
$ws.Range("J129").Value = 1933.3334
$ws.Range("L129").Value = 5800.0002
$ws.Range("N129").Value = -15800.0002

$ws.Range("H139").Value = 1864.6923
$ws.Range("I139").Value = 1570.0834
$ws.Range("J139").Value = 5400
$ws.Range("K139").Value = 4710.2502
$ws.Range("L139").Value = 16200
$ws.Range("M139").Value = 429.7497999999996
$ws.Range("N139").Value = -26480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 19999.5
$ws.Range("I18").Value = 19999.5
$ws.Range("K18").Value = 19999.5
$ws.Range("M18").Value = -19706.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5334.875
$ws.Range("I7").Value = 4136
$ws.Range("K7").Value = 4136
$ws.Range("M7").Value = -4024

$ws.Range("H22").Value = 1794.2
$ws.Range("I22").Value = 2042.75
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 2042.75
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -1747.75
$ws.Range("N22").Value = -1390

$ws.Range("H27").Value = 1794.2
$ws.Range("I27").Value = 2042.75
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 2042.75
$ws.Range("L27").Value = 800
$ws.Range("M27").Value = -1935.75
$ws.Range("N27").Value = -1014

$ws.Range("H126").Value = 5334.875
$ws.Range("I126").Value = 4136
$ws.Range("K126").Value = 12408
$ws.Range("M126").Value = -9938

$ws.Range("H136").Value = 1922.9
$ws.Range("I136").Value = 1713.3448
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 5140.0344
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -2590.0344
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 64000
$ws.Range("I25").Value = 50000
$ws.Range("K25").Value = 50000
$ws.Range("M25").Value = -49707

$ws.Range("H100").Value = 71429290
$ws.Range("I100").Value = 810.375
$ws.Range("J100").Value = 166667250
$ws.Range("K100").Value = 1620.75
$ws.Range("L100").Value = 333334500
$ws.Range("M100").Value = -1079.75
$ws.Range("N100").Value = -333335582

$ws.Range("H136").Value = 4056.724
$ws.Range("I136").Value = 2415.7144
$ws.Range("K136").Value = 7247.1432
$ws.Range("M136").Value = -4697.1432
